$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 145.7087706666667
$ws.Range("N2").Value = 437.126312
$ws.Range("O2").Value = 0.5445232453600627
$ws.Range("P2").Value = 0.5461141113270247
$ws.Range("Q2").Value = 2436.099011285993
$ws.Range("R2").Value = 21924.89110157394
$ws.Range("S2").Value = 0.01585730503782149
$ws.Range("T2").Value = 0.01670710142172703
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.003349722554576428
$ws.Range("P3").Value = 0.003359509023117945
$ws.Range("Q3").Value = 14.98605591739267
$ws.Range("R3").Value = 134.874503256534
$ws.Range("S3").Value = 0.00009754876911611989
$ws.Range("T3").Value = 0.000102776428611324
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 51.59199533333333
$ws.Range("N4").Value = 154.775986
$ws.Range("O4").Value = 0.1928026748491032
$ws.Range("P4").Value = 0.1933659624890163
$ws.Range("Q4").Value = 862.5644719035231
$ws.Range("R4").Value = 7763.080247131707
$ws.Range("S4").Value = 0.005614692950653559
$ws.Range("T4").Value = 0.005915585552191157
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 2.338518
$ws.Range("N5").Value = 4.677036
$ws.Range("O5").Value = 0.0087391953474509
$ws.Range("P5").Value = 0.005843151713055659
$ws.Range("Q5").Value = 39.09758734226801
$ws.Range("R5").Value = 234.585524053608
$ws.Range("S5").Value = 0.0002544980174684811
$ws.Range("T5").Value = 0.0001787577472688684
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("M6").Value = 67.05398933333333
$ws.Range("N6").Value = 201.161968
$ws.Range("O6").Value = 0.2505851618888069
$ws.Range("P6").Value = 0.2513172654477853
$ws.Range("Q6").Value = 1121.072920801767
$ws.Range("R6").Value = 10089.6562872159
$ws.Range("S6").Value = 0.007297402606559372
$ws.Range("T6").Value = 0.007688471980085721
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 145.7087706666667
$ws.Range("N7").Value = 437.126312
$ws.Range("O7").Value = 0.5445232453600627
$ws.Range("P7").Value = 0.5461141113270247
$ws.Range("Q7").Value = 22657.77051937846
$ws.Range("R7").Value = 203919.9346744061
$ws.Range("S7").Value = 0.1474862790626383
$ws.Range("T7").Value = 0.1553901004449085
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.003349722554576428
$ws.Range("P8").Value = 0.003359509023117945
$ws.Range("S8").Value = 0.0009072856295418445
$ws.Range("T8").Value = 0.000955907261358549
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 51.59199533333333
$ws.Range("N9").Value = 154.775986
$ws.Range("O9").Value = 0.1928026748491032
$ws.Range("P9").Value = 0.1933659624890163
$ws.Range("Q9").Value = 8022.575343619518
$ws.Range("R9").Value = 72203.17809257566
$ws.Range("S9").Value = 0.05222136859926885
$ws.Range("T9").Value = 0.05501992296222095
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 2.338518
$ws.Range("N10").Value = 4.677036
$ws.Range("O10").Value = 0.0087391953474509
$ws.Range("P10").Value = 0.005843151713055659
$ws.Range("Q10").Value = 363.6404586835021
$ws.Range("R10").Value = 2181.842752101012
$ws.Range("S10").Value = 0.002367045695073621
$ws.Range("T10").Value = 0.0016625974549536
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("M11").Value = 67.05398933333333
$ws.Range("N11").Value = 201.161968
$ws.Range("O11").Value = 0.2505851618888069
$ws.Range("P11").Value = 0.2513172654477853
$ws.Range("Q11").Value = 10426.92142533519
$ws.Range("R11").Value = 93842.29282801665
$ws.Range("S11").Value = 0.06787198421777346
$ws.Range("T11").Value = 0.07150925843424287
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 145.7087706666667
$ws.Range("N12").Value = 437.126312
$ws.Range("O12").Value = 0.5445232453600627
$ws.Range("P12").Value = 0.5461141113270247
$ws.Range("Q12").Value = 28392.05968054592
$ws.Range("R12").Value = 255528.5371249134
$ws.Range("S12").Value = 0.1848125010193168
$ws.Range("T12").Value = 0.1947166426557529
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.003349722554576428
$ws.Range("P13").Value = 0.003359509023117945
$ws.Range("Q13").Value = 174.6583336766766
$ws.Range("R13").Value = 1571.92500309009
$ws.Range("S13").Value = 0.001136903903198344
$ws.Range("T13").Value = 0.001197830827633594
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 51.59199533333333
$ws.Range("N14").Value = 154.775986
$ws.Range("O14").Value = 0.1928026748491032
$ws.Range("P14").Value = 0.1933659624890163
$ws.Range("Q14").Value = 10052.95016793073
$ws.Range("R14").Value = 90476.55151137657
$ws.Range("S14").Value = 0.06543769223022833
$ws.Range("T14").Value = 0.06894451221607961
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 2.338518
$ws.Range("N15").Value = 4.677036
$ws.Range("O15").Value = 0.0087391953474509
$ws.Range("P15").Value = 0.005843151713055659
$ws.Range("Q15").Value = 455.67155852218
$ws.Range("R15").Value = 2734.02935113308
$ws.Range("S15").Value = 0.002966103950237781
$ws.Range("T15").Value = 0.00208337206546398
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("M16").Value = 67.05398933333333
$ws.Range("N16").Value = 201.161968
$ws.Range("O16").Value = 0.2505851618888069
$ws.Range("P16").Value = 0.2513172654477853
$ws.Range("Q16").Value = 13065.79458642167
$ws.Range("R16").Value = 117592.151277795
$ws.Range("S16").Value = 0.08504920750697749
$ws.Range("T16").Value = 0.08960701281002736
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 145.7087706666667
$ws.Range("N17").Value = 437.126312
$ws.Range("O17").Value = 0.5445232453600627
$ws.Range("P17").Value = 0.5461141113270247
$ws.Range("Q17").Value = 12068.98491135221
$ws.Range("R17").Value = 72413.90946811324
$ws.Range("S17").Value = 0.07856067193884225
$ws.Range("T17").Value = 0.05518050348449288
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.003349722554576428
$ws.Range("P18").Value = 0.003359509023117945
$ws.Range("Q18").Value = 74.24430694720199
$ws.Range("R18").Value = 445.4658416832119
$ws.Range("S18").Value = 0.0004832786422592645
$ws.Range("T18").Value = 0.000339451765686633
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 51.59199533333333
$ws.Range("N19").Value = 154.775986
$ws.Range("O19").Value = 0.1928026748491032
$ws.Range("P19").Value = 0.1933659624890163
$ws.Range("Q19").Value = 4273.339280646323
$ws.Range("R19").Value = 25640.03568387794
$ws.Range("S19").Value = 0.02781645745488055
$ws.Range("T19").Value = 0.01953809825748678
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 2.338518
$ws.Range("N20").Value = 4.677036
$ws.Range("O20").Value = 0.0087391953474509
$ws.Range("P20").Value = 0.005843151713055659
$ws.Range("Q20").Value = 193.698281358036
$ws.Range("R20").Value = 774.7931254321439
$ws.Range("S20").Value = 0.001260840679531624
$ws.Range("T20").Value = 0.0005904041788614607
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("M21").Value = 67.05398933333333
$ws.Range("N21").Value = 201.161968
$ws.Range("O21").Value = 0.2505851618888069
$ws.Range("P21").Value = 0.2513172654477853
$ws.Range("Q21").Value = 5554.048543593311
$ws.Range("R21").Value = 33324.29126155987
$ws.Range("S21").Value = 0.0361529812797448
$ws.Range("T21").Value = 0.02539361820930937
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 145.7087706666667
$ws.Range("N22").Value = 437.126312
$ws.Range("O22").Value = 0.5445232453600627
$ws.Range("P22").Value = 0.5461141113270247
$ws.Range("Q22").Value = 18098.17424775031
$ws.Range("R22").Value = 162883.5682297528
$ws.Range("S22").Value = 0.1178064883014438
$ws.Range("T22").Value = 0.1241197633201434
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.003349722554576428
$ws.Range("P23").Value = 0.003359509023117945
$ws.Range("Q23").Value = 111.3338374273743
$ws.Range("R23").Value = 1002.004536846369
$ws.Range("S23").Value = 0.0007247056104608552
$ws.Range("T23").Value = 0.0007635427398278454
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 51.59199533333333
$ws.Range("N24").Value = 154.775986
$ws.Range("O24").Value = 0.1928026748491032
$ws.Range("P24").Value = 0.1933659624890163
$ws.Range("Q24").Value = 6408.131213101997
$ws.Range("R24").Value = 57673.18091791797
$ws.Range("S24").Value = 0.04171246361407186
$ws.Range("T24").Value = 0.04394784350103782
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 2.338518
$ws.Range("N25").Value = 4.677036
$ws.Range("O25").Value = 0.0087391953474509
$ws.Range("P25").Value = 0.005843151713055659
$ws.Range("Q25").Value = 290.462310894938
$ws.Range("R25").Value = 1742.773865369628
$ws.Range("S25").Value = 0.001890707005139391
$ws.Range("T25").Value = 0.001328020266507751
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("M26").Value = 67.05398933333333
$ws.Range("N26").Value = 201.161968
$ws.Range("O26").Value = 0.2505851618888069
$ws.Range("P26").Value = 0.2513172654477853
$ws.Range("Q26").Value = 8328.632363096851
$ws.Range("R26").Value = 74957.69126787166
$ws.Range("S26").Value = 0.0542135862777517
$ws.Range("T26").Value = 0.08960701281002736
